# Update on 7th June 2020 - refresh the last few rows of daily testing
# data and append the new day's row (row 68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (revised figures for 65-67) ---
$ws.Range("J65").Value  = 230145
$ws.Range("AK65").Value = 4685525

$ws.Range("J66").Value  = 236506
$ws.Range("AK66").Value = 4837463

$ws.Range("D67").Value  = 10790
$ws.Range("J67").Value  = 236506
$ws.Range("AK67").Value = 4985013

# --- New row 68 (states_tested_data for 2020-06-06, serial 43988) ---
$row68 = @(
    43988, 6677, 436335, 11261, 146605, 95473, 5237, 89392, 14780, 246873,
    27402, 245606, 141688, 44509, 211880, 84444, 372582, 107796, 10164,
    200913, 538009, 14629, 10066, 1991, 4061, 172598, 7963, 115974, 494480,
    5005, 576695, 23388, 33331, 355085, 36638, 261288, 5150818
)

for ($i = 0; $i -lt $row68.Length; $i++) {
    $ws.Cells.Item(68, $i + 1).Value = $row68[$i]
}
